# Implement File Upload with AutoIT
# Update the SkillDetails sheet: change the "From Date" values in H2/H3
# from 44109 (10/5/2020) to 44119 (10/15/2020), and move the active
# selection from P4 to H3.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H2").Value = 44119
$ws.Range("H3").Value = 44119

$ws.Range("H3").Select()
